$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "61.590.15"
Set-TextValue $ws.Range("E2") "  -2.82%  "
Set-TextValue $ws.Range("D3") "3.389.61"
Set-TextValue $ws.Range("E3") "  -2.64%  "
Set-TextValue $ws.Range("E4") "  -0.07%  "
Set-TextValue $ws.Range("D5") "405.38"
Set-TextValue $ws.Range("E5") "  -2.82%  "
Set-TextValue $ws.Range("D6") "133.58"
Set-TextValue $ws.Range("E6") "  +8.48%  "
Set-TextValue $ws.Range("E7") "  -2.03%  "
Set-TextValue $ws.Range("E8") "  -0.09%  "
Set-TextValue $ws.Range("E9") "  -2.28%  "
Set-TextValue $ws.Range("E10") "  -9.75%  "
Set-TextValue $ws.Range("D11") "42.49"
Set-TextValue $ws.Range("E11") "  +2.22%  "
Set-TextValue $ws.Range("E12") "  -1.16%  "
Set-TextValue $ws.Range("D13") "3.908.61"
Set-TextValue $ws.Range("E13") "  -3.20%  "
Set-TextValue $ws.Range("D14") "8.41"
Set-TextValue $ws.Range("E14") "  -2.49%  "
Set-TextValue $ws.Range("D15") "19.78"
Set-TextValue $ws.Range("E15") "  -1.08%  "
Set-TextValue $ws.Range("D16") "3.389.04"
Set-TextValue $ws.Range("E16") "  -2.72%  "
Set-TextValue $ws.Range("D17") "61.555.18"
Set-TextValue $ws.Range("E17") "  -2.76%  "
Set-TextValue $ws.Range("D19") "11.06"
Set-TextValue $ws.Range("E19") "  -0.05%  "
Set-TextValue $ws.Range("E20") "  -11.39%  "
Set-TextValue $ws.Range("E21") "  -3.43%  "
Set-TextValue $ws.Range("D22") "85.54"
Set-TextValue $ws.Range("E22") "  +3.39%  "
Set-TextValue $ws.Range("D23") "315.83"
Set-TextValue $ws.Range("E23") "  -0.74%  "
Set-TextValue $ws.Range("D24") "12.79"
Set-TextValue $ws.Range("E24") "  -0.91%  "
Set-TextValue $ws.Range("D25") "3.12"
Set-TextValue $ws.Range("E25") "  -1.91%  "
Set-TextValue $ws.Range("E26") "  +11.03%  "
Set-TextValue $ws.Range("D27") "29.54"
Set-TextValue $ws.Range("E27") "  -5.40%  "
Set-TextValue $ws.Range("D28") "8.30"
Set-TextValue $ws.Range("E28") "  +5.12%  "
Set-TextValue $ws.Range("E29") "  -2.15%  "
Set-TextValue $ws.Range("E30") "  +0.15%  "
Set-TextValue $ws.Range("E31") "  -2.65%  "
Set-TextValue $ws.Range("D32") "2.64"
Set-TextValue $ws.Range("E32") "  +3.23%  "
Set-TextValue $ws.Range("D33") "11.34"
Set-TextValue $ws.Range("E33") "  -2.45%  "
Set-TextValue $ws.Range("E34") "  -0.71%  "
Set-TextValue $ws.Range("D35") "41.12"
Set-TextValue $ws.Range("E35") "  -2.19%  "
Set-TextValue $ws.Range("D36") "0.0479"
Set-TextValue $ws.Range("E36") "  -2.01%  "
Set-TextValue $ws.Range("D37") "51.73"
Set-TextValue $ws.Range("E37") "  -0.71%  "
Set-TextValue $ws.Range("E38") "  +0.08%  "
Set-TextValue $ws.Range("D39") "3.42"
Set-TextValue $ws.Range("E39") "  -1.57%  "
Set-TextValue $ws.Range("E40") "  -3.78%  "
Set-TextValue $ws.Range("D41") "139.54"
Set-TextValue $ws.Range("E41") "  +2.82%  "
Set-TextValue $ws.Range("E42") "  -1.55%  "
Set-TextValue $ws.Range("E43") "  -1.60%  "
Set-TextValue $ws.Range("E44") "  +4.40%  "
Set-TextValue $ws.Range("D45") "3.98"
Set-TextValue $ws.Range("E45") "  +1.87%  "
Set-TextValue $ws.Range("D46") "16.59"
Set-TextValue $ws.Range("E46") "  -2.19%  "
Set-TextValue $ws.Range("E47") "  -2.10%  "
Set-TextValue $ws.Range("D48") "21.32"
Set-TextValue $ws.Range("E48") "  -3.13%  "
Set-TextValue $ws.Range("D49") "2.118.93"
Set-TextValue $ws.Range("E49") "  -3.18%  "
Set-TextValue $ws.Range("E50") "  -6.03%  "
Set-TextValue $ws.Range("E51") "  -0.26%  "
